$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.357.12'
$ws.Range('E2').Value = '  +4.38%  '
$ws.Range('D3').Value = '1.736.39'
$ws.Range('E3').Value = '  +3.27%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '''220.25'
$ws.Range('E5').Value = '  +2.27%  '
$ws.Range('D6').Value = '''0.523'
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = '''24.28'
$ws.Range('E8').Value = '  +12.76%  '
$ws.Range('E9').Value = '  +4.13%  '
$ws.Range('D10').Value = '''0.0637'
$ws.Range('E10').Value = '  +2.05%  '
$ws.Range('D11').Value = '''0.0896'
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('D12').Value = '1.982.50'
$ws.Range('E12').Value = '  +3.34%  '
$ws.Range('D13').Value = '1.733.71'
$ws.Range('E13').Value = '  +3.28%  '
$ws.Range('E14').Value = '  +3.35%  '
$ws.Range('E15').Value = '  +4.59%  '
$ws.Range('D16').Value = '''67.80'
$ws.Range('E16').Value = '  +2.11%  '
$ws.Range('D17').Value = '28.346.92'
$ws.Range('E17').Value = '  +4.40%  '
$ws.Range('D18').Value = '''243.43'
$ws.Range('E18').Value = '  +1.94%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0758'
$ws.Range('E19').Value = '  +2.08%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '''8.02'
$ws.Range('E20').Value = '  -0.69%  '
$ws.Range('D21').Value = '''0.999'
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('E22').Value = '  +2.90%  '
$ws.Range('E23').Value = '  +2.62%  '
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('D25').Value = '''149.76'
$ws.Range('E25').Value = '  +0.99%  '
$ws.Range('E26').Value = '  +4.03%  '
$ws.Range('D27').Value = '''16.76'
$ws.Range('E27').Value = '  +2.83%  '
$ws.Range('E28').Value = '  +1.57%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('E30').Value = '  +3.24%  '
$ws.Range('E31').Value = '  +3.24%  '
$ws.Range('D32').Value = '''3.44'
$ws.Range('E32').Value = '  +2.04%  '
$ws.Range('D33').Value = '1.505.42'
$ws.Range('E33').Value = '  -4.08%  '
$ws.Range('E34').Value = '  +2.02%  '
$ws.Range('E35').Value = '  -1.51%  '
$ws.Range('D36').Value = '''0.970'
$ws.Range('E36').Value = '  +3.66%  '
$ws.Range('D37').Value = '''0.607'
$ws.Range('E37').Value = '  +0.81%  '
$ws.Range('E38').Value = '  +0.49%  '
$ws.Range('E39').Value = '  +1.51%  '
$ws.Range('E40').Value = '  +1.64%  '
$ws.Range('D41').Value = '''70.86'
$ws.Range('E41').Value = '  +2.47%  '
$ws.Range('D42').Value = '''5.73'
$ws.Range('E42').Value = '  +2.30%  '
$ws.Range('D43').Value = '''0.999'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('E44').Value = '  +2.24%  '
$ws.Range('D45').Value = '1.884.93'
$ws.Range('D46').Value = '''0.802'
$ws.Range('E46').Value = '  +1.95%  '
$ws.Range('E47').Value = '  +9.66%  '
$ws.Range('D48').Value = '''91.21'
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('D49').Value = '0.0₆0111'
$ws.Range('E49').Value = '  +3.59%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '''0.105'
$ws.Range('E50').Value = '  +1.27%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''8.23'
$ws.Range('E51').Value = '  +1.08%  '

Write-Host "Applied crypto list update"
